# This document shuffles several fixed paragraph/run "slots" (labels, headings,
# etc. keep their position/formatting) so that the VALUE text that used to sit in
# one slot ends up in a different slot (see commit diff). Because several values
# trade places in a cycle, a direct "search old, replace new" pass would clobber
# data that a later step still needs to find. So we do it in two passes:
#   Phase 1: stamp every value that is about to move with a unique placeholder
#            token (so the original text disappears from the document).
#   Phase 2: turn each placeholder into the text that should finally occupy that
#            slot.
# Multi-line values (paragraph programs / bibliography) are written/matched as one
# logical Find target using [char]11 (Word's manual-line-break / <w:br/> marker in
# Range.Text) to join the lines, so formatting (w:br elements) is preserved intact.

$d = $word.ActiveDocument

# ---- Phase 1: tag each value that needs to move with a unique placeholder ----
$null = $d.Content.Find.Execute("Apresentar aos alunos os conceitos introdutórios de Física e em particular, da Mecânica incluindo cinemática e dinâmica, além de conceitos de estatística básica e análise de dados.", $true, $false, $false, $false, $false, $true, 1, $false, "@@A@@", 2)
$null = $d.Content.Find.Execute("Presenting to the students the introductory concepts of Physics and in particular, of Mechanics including kinematics and dynamics, including basic concepts of statistical and data analysis.", $true, $false, $false, $false, $false, $true, 1, $false, "@@B@@", 2)
$null = $d.Content.Find.Execute("Introdução a física, Cinemática, Dinâmica, Energia, Momento linear, Rotação.", $true, $false, $false, $false, $false, $true, 1, $false, "@@C@@", 2)
$null = $d.Content.Find.Execute("Introduction to Physics, Kinematics, Dynamics, Energy, Linear momentum, Rotation", $true, $false, $false, $false, $false, $true, 1, $false, "@@D@@", 2)
$null = $d.Content.Find.Execute("1) Introdução a Física: sistemas de unidades, revisão de vetores, análise dimensional." + [char]11 + "2) Cinemática: movimento unidimensional, queda livre, movimento bidimensional, projéteis. " + [char]11 + "3) Dinâmica: leis de Newton, forças, força de atrito, força de resistência do ar, velocidade terminal, movimento circular uniforme, gravitação, aplicações." + [char]11 + "4) Energia: trabalho, forças conservativas, conservação de energia mecânica, atrito, aplicações." + [char]11 + "5)  Momento linear: centro de massa, sistema de partículas, conservação do momento linear, colisões, impulso." + [char]11 + "6) Rotação: variáveis do movimento rotacional, energia cinética rotacional, momento de inércia, torque, rolamento, conservação do momento angular.", $true, $false, $false, $false, $false, $true, 1, $false, "@@E@@", 2)
$null = $d.Content.Find.Execute("As avaliações serão compostas por provas, projetos, seminários e outras formas que serão utilizadas para a composição das notas. A média final (NF) é calculada pela média simples das notas (N), levando em conta o número n de avaliações, sendo no mínimo duas avaliações: NF= (N1+...+Nn)/n.", $true, $false, $false, $false, $false, $true, 1, $false, "@@G@@", 2)
$null = $d.Content.Find.Execute("NF ≥ 5,0", $true, $false, $false, $false, $false, $true, 1, $false, "@@H@@", 2)
$null = $d.Content.Find.Execute("(NF+REC)/2 ≥ 5,0, onde REC é uma prova de recuperação a ser aplicada, seguindo as regras da EEL.", $true, $false, $false, $false, $false, $true, 1, $false, "@@I@@", 2)
$null = $d.Content.Find.Execute("HALLIDAY, D; RESNICK, R. Fundamentos de Física. Vol.1, LTC (2008)." + [char]11 + "SEARS, F. W.; ZEMANSKY, M. W.; YOUNG, H. D.; FREEDMAN, R. A. Física I, Vol. 1, Pearson Addison Wesley (2009)." + [char]11 + "JEWETT Jr, John W.; SERWAY, Raymond A. Princípios de Física. Vol. 1, Thomson Pioneira (2008)." + [char]11 + "NUSSENZVEIG, H.M. Curso de Física Básica. Vol. 1, Edgard Blucher (2008)." + [char]11 + "TIPLER, P.; MOSCA, G. Física para Cientistas e Engenheiros. Vol.1, LTC (2008).", $true, $false, $false, $false, $false, $true, 1, $false, "@@J@@", 2)
$null = $d.Content.Find.Execute("8711686 - Flavia Reis Cardoso Rojas", $true, $false, $false, $false, $false, $true, 1, $false, "@@K@@", 2)

# ---- Phase 2: resolve each placeholder to the value that now belongs there ----
$null = $d.Content.Find.Execute("@@A@@", $true, $false, $false, $false, $false, $true, 1, $false, "Introdução a física, Cinemática, Dinâmica, Energia, Momento linear, Rotação.", 2)
$null = $d.Content.Find.Execute("@@B@@", $true, $false, $false, $false, $false, $true, 1, $false, "Introduction to Physics, Kinematics, Dynamics, Energy, Linear momentum, Rotation", 2)
$null = $d.Content.Find.Execute("@@K@@", $true, $false, $false, $false, $false, $true, 1, $false, "Apresentar aos alunos os conceitos introdutórios de Física e em particular, da Mecânica incluindo cinemática e dinâmica, além de conceitos de estatística básica e análise de dados.", 2)
$null = $d.Content.Find.Execute("@@C@@", $true, $false, $false, $false, $false, $true, 1, $false, "1) Introdução a Física: sistemas de unidades, revisão de vetores, análise dimensional." + [char]11 + "2) Cinemática: movimento unidimensional, queda livre, movimento bidimensional, projéteis. " + [char]11 + "3) Dinâmica: leis de Newton, forças, força de atrito, força de resistência do ar, velocidade terminal, movimento circular uniforme, gravitação, aplicações." + [char]11 + "4) Energia: trabalho, forças conservativas, conservação de energia mecânica, atrito, aplicações." + [char]11 + "5)  Momento linear: centro de massa, sistema de partículas, conservação do momento linear, colisões, impulso." + [char]11 + "6) Rotação: variáveis do movimento rotacional, energia cinética rotacional, momento de inércia, torque, rolamento, conservação do momento angular.", 2)
$null = $d.Content.Find.Execute("@@D@@", $true, $false, $false, $false, $false, $true, 1, $false, "Presenting to the students the introductory concepts of Physics and in particular, of Mechanics including kinematics and dynamics, including basic concepts of statistical and data analysis.", 2)
$null = $d.Content.Find.Execute("@@E@@", $true, $false, $false, $false, $false, $true, 1, $false, "As avaliações serão compostas por provas, projetos, seminários e outras formas que serão utilizadas para a composição das notas. A média final (NF) é calculada pela média simples das notas (N), levando em conta o número n de avaliações, sendo no mínimo duas avaliações: NF= (N1+...+Nn)/n.", 2)
$null = $d.Content.Find.Execute("@@G@@", $true, $false, $false, $false, $false, $true, 1, $false, "NF ≥ 5,0", 2)
$null = $d.Content.Find.Execute("@@H@@", $true, $false, $false, $false, $false, $true, 1, $false, "(NF+REC)/2 ≥ 5,0, onde REC é uma prova de recuperação a ser aplicada, seguindo as regras da EEL.", 2)
$null = $d.Content.Find.Execute("@@I@@", $true, $false, $false, $false, $false, $true, 1, $false, "HALLIDAY, D; RESNICK, R. Fundamentos de Física. Vol.1, LTC (2008)." + [char]11 + "SEARS, F. W.; ZEMANSKY, M. W.; YOUNG, H. D.; FREEDMAN, R. A. Física I, Vol. 1, Pearson Addison Wesley (2009)." + [char]11 + "JEWETT Jr, John W.; SERWAY, Raymond A. Princípios de Física. Vol. 1, Thomson Pioneira (2008)." + [char]11 + "NUSSENZVEIG, H.M. Curso de Física Básica. Vol. 1, Edgard Blucher (2008)." + [char]11 + "TIPLER, P.; MOSCA, G. Física para Cientistas e Engenheiros. Vol.1, LTC (2008).", 2)
$null = $d.Content.Find.Execute("@@J@@", $true, $false, $false, $false, $false, $true, 1, $false, "8711686 - Flavia Reis Cardoso Rojas", 2)
